# Added periodic & upfront related scenarios
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Update the "repaymentstrategy" value (B17) from "RBI (India)" to the new scenario value
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Reflect the author's new scroll position / selection on that sheet
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("B17").Select()
